# Add data for 2022-08-24:
#  - roll the "through" date in the sheet name and header label from
#    August 15 to August 16
#  - update the affected neighborhood/month cell counts

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet tab and update the matching header label in B1.
$wb.Worksheets.Item(1).Name = "Through 2022-08-16"
$ws.Range("B1").Value = "August 2022 (through August 16)"

# Map of cell address -> new value for the updated counts.
$updates = @{
    "B2"   = 12
    "Z2"   = 2
    "AX3"  = 3
    "AH4"  = 3
    "AX4"  = 1
    "B5"   = 6
    "R5"   = 9
    "AH5"  = 2
    "AP5"  = 1
    "AX9"  = 6
    "J13"  = 6
    "J14"  = 1
    "AX14" = 1
    "AP15" = 1
    "B16"  = 1
    "AP17" = 4
    "R19"  = 1
    "J21"  = 1
    "R21"  = 2
    "AX23" = 1
    "B26"  = 2
    "R26"  = 1
    "AP33" = 1
    "B37"  = 4
    "J37"  = 2
    "BF48" = 1
    "BF54" = 1
    "R55"  = 1
    "R60"  = 2
    "BF63" = 1
    "R66"  = 2
    "BF89" = 1
    "AP90" = 1
    "R93"  = 1
    "AP96" = 1
    "AX97" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
